# Generate Report for Handback
# Updates the localization-status report to reflect that handback has
# completed and is in sync with en-US: statuses, handback datetimes,
# target/handback file columns (with hyperlinks), and related column widths.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Overview sheet: flip the per-locale status cells from "Ready for handoff"
# to "Handed back: in sync with en-US"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"

# widen the status columns to fit the longer text
$wsOverview.Columns.Item(5).ColumnWidth = 29.15
$wsOverview.Columns.Item(6).ColumnWidth = 29.15

# ---------------------------------------------------------------------------
# Helper data shared by both locale sheets
# ---------------------------------------------------------------------------
$file1Name = "1e4daa63-87a3-4c3b-be57-3f2ff07dceb2.md"
$file1Url  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e3611e5106f5391dafbbef73fea35cf0eab388ac/e2e/1e4daa63-87a3-4c3b-be57-3f2ff07dceb2.md"
$file2Name = "d89496bf-a4f9-4407-a687-138de49aa108.md"
$file2Url  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e3611e5106f5391dafbbef73fea35cf0eab388ac/e2e/d89496bf-a4f9-4407-a687-138de49aa108.md"

function Set-HandbackColumns($ws, $xlf1, $xlf2, $handbackDate) {
    # Latest Target File (I) - link back to the source markdown file
    $ws.Range("I2").Value = $file1Name
    $ws.Range("I3").Value = $file2Name

    $ws.Hyperlinks.Add($ws.Range("I2"), $file1Url, "", "", $file1Name)
    $ws.Hyperlinks.Add($ws.Range("I3"), $file2Url, "", "", $file2Name)

    # match the hyperlink look already used for column A
    $ws.Range("I2").Font.Underline = 2
    $ws.Range("I2").Font.Color = 15570276
    $ws.Range("I3").Font.Underline = 2
    $ws.Range("I3").Font.Color = 15570276

    # Latest Handback File (J) - the xliff that was last handed off/back
    $ws.Range("J2").Value = $xlf1
    $ws.Range("J3").Value = $xlf2

    # Latest Handback DateTime (K)
    $ws.Range("K2").Value = $handbackDate
    $ws.Range("K3").Value = $handbackDate

    # widen columns to fit the new content
    $ws.Columns.Item(3).ColumnWidth = 29.15
    $ws.Columns.Item(9).ColumnWidth = 39.15
    $ws.Columns.Item(10).ColumnWidth = 39.15
}

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
Set-HandbackColumns $wsZhCn `
    "1e4daa63-87a3-4c3b-be57-3f2ff07dceb2.df79be2adc3d8825dffc00a54a89d9846acc4a34.zh-cn.xlf" `
    "d89496bf-a4f9-4407-a687-138de49aa108.e2b2914fa0e0c81c501573a3ce74dfcd4d01df49.zh-cn.xlf" `
    "2016-10-14 08:45:41"

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
Set-HandbackColumns $wsDeDe `
    "1e4daa63-87a3-4c3b-be57-3f2ff07dceb2.df79be2adc3d8825dffc00a54a89d9846acc4a34.de-de.xlf" `
    "d89496bf-a4f9-4407-a687-138de49aa108.e2b2914fa0e0c81c501573a3ce74dfcd4d01df49.de-de.xlf" `
    "2016-10-14 08:45:57"

Write-Output "Handback report generated."
